$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.941.31"
$ws.Range("E2").Value = "  -3.30%  "
$ws.Range("D3").Value = "3.280.51"
$ws.Range("E3").Value = "  -4.18%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "'554.40"
$ws.Range("E5").Value = "  -4.20%  "
$ws.Range("D6").Value = "'140.17"
$ws.Range("E6").Value = "  -8.50%  "
$ws.Range("E7").Value = "  -0.09%  "
$ws.Range("D8").Value = "3.282.64"
$ws.Range("E8").Value = "  -4.12%  "
$ws.Range("D9").Value = "'0.464"
$ws.Range("E9").Value = "  -4.09%  "
$ws.Range("D10").Value = "'7.85"
$ws.Range("E10").Value = "  -2.64%  "
$ws.Range("E11").Value = "  -5.91%  "
$ws.Range("E12").Value = "  -3.75%  "
$ws.Range("D13").Value = "3.839.19"
$ws.Range("E13").Value = "  -4.23%  "
$ws.Range("E14").Value = "  +0.01%  "
$ws.Range("D15").Value = "'26.54"
$ws.Range("E15").Value = "  -6.88%  "
$ws.Range("D16").Value = "3.271.77"
$ws.Range("E16").Value = "  -4.36%  "
$ws.Range("E17").Value = "  -5.45%  "
$ws.Range("D18").Value = "59.981.89"
$ws.Range("E18").Value = "  -3.28%  "
$ws.Range("D19").Value = "'6.01"
$ws.Range("E19").Value = "  -8.32%  "
$ws.Range("D20").Value = "'13.63"
$ws.Range("E20").Value = "  -5.93%  "
$ws.Range("D21").Value = "'8.46"
$ws.Range("E21").Value = "  -5.57%  "
$ws.Range("D22").Value = "'371.02"
$ws.Range("E22").Value = "  -2.90%  "
$ws.Range("D24").Value = "'72.14"
$ws.Range("E24").Value = "  -4.14%  "
$ws.Range("D25").Value = "'0.529"
$ws.Range("E25").Value = "  -7.70%  "
$ws.Range("D26").Value = "3.406.48"
$ws.Range("E26").Value = "  -4.30%  "
$ws.Range("E27").Value = "  -9.47%  "
$ws.Range("E28").Value = "  -4.33%  "
$ws.Range("D29").Value = "'0.997"
$ws.Range("E29").Value = "  -0.17%  "
$ws.Range("D30").Value = "'7.01"
$ws.Range("E30").Value = "  -8.64%  "
$ws.Range("E31").Value = "  +0.00%  "
$ws.Range("D32").Value = "'2.01"
$ws.Range("E32").Value = "  -5.57%  "
$ws.Range("D33").Value = "'7.40"
$ws.Range("E33").Value = "  -6.10%  "
$ws.Range("D34").Value = "'22.45"
$ws.Range("E35").Value = "  -8.21%  "
$ws.Range("D36").Value = "'166.19"
$ws.Range("E36").Value = "  -1.60%  "
$ws.Range("E37").Value = "  -9.25%  "
$ws.Range("E38").Value = "  -5.68%  "
$ws.Range("D39").Value = "'6.58"
$ws.Range("E39").Value = "  -5.57%  "
$ws.Range("D40").Value = "3.309.17"
$ws.Range("E40").Value = "  -4.31%  "
$ws.Range("E41").Value = "  -8.41%  "
$ws.Range("D42").Value = "'25.54"
$ws.Range("E42").Value = "  -18.22%  "
$ws.Range("D43").Value = "'41.44"
$ws.Range("E43").Value = "  -2.87%  "
$ws.Range("D44").Value = "'0.743"
$ws.Range("E44").Value = "  -4.74%  "
$ws.Range("E45").Value = "  -8.18%  "
$ws.Range("E46").Value = "  -4.26%  "
$ws.Range("D47").Value = "'1.56"
$ws.Range("E47").Value = "  -7.36%  "
$ws.Range("E48").Value = "  -0.03%  "
$ws.Range("D49").Value = "2.306.36"
$ws.Range("E49").Value = "  -9.53%  "
$ws.Range("D50").Value = "'6.30"
$ws.Range("E50").Value = "  -7.83%  "
$ws.Range("D51").Value = "'21.30"
$ws.Range("E51").Value = "  -5.43%  "